$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: new date columns DO1 (2026-02-11) and DP1 (2026-02-12) ---
$ws.Range("DO1").Value = 46064
$ws.Range("DN1").Copy()
$ws.Range("DO1").PasteSpecial(-4122)

$ws.Range("DP1").Value = 46065
$ws.Range("DN1").Copy()
$ws.Range("DP1").PasteSpecial(-4122)

# --- Data rows: attendance values for the two new training dates ---
# Row 2
$ws.Range("DO2").Value = "P"
$ws.Range("DN2").Copy()
$ws.Range("DO2").PasteSpecial(-4122)
$ws.Range("DP2").Value = "P"
$ws.Range("DN2").Copy()
$ws.Range("DP2").PasteSpecial(-4122)

# Row 3
$ws.Range("DO3").Value = "P"
$ws.Range("DN3").Copy()
$ws.Range("DO3").PasteSpecial(-4122)
$ws.Range("DP3").Value = "M"
$ws.Range("DN3").Copy()
$ws.Range("DP3").PasteSpecial(-4122)

# Row 4
$ws.Range("DO4").Value = "P"
$ws.Range("DN4").Copy()
$ws.Range("DO4").PasteSpecial(-4122)
$ws.Range("DP4").Value = "P"
$ws.Range("DN4").Copy()
$ws.Range("DP4").PasteSpecial(-4122)

# Row 5
$ws.Range("DO5").Value = "P"
$ws.Range("DN5").Copy()
$ws.Range("DO5").PasteSpecial(-4122)
$ws.Range("DP5").Value = "P"
$ws.Range("DN5").Copy()
$ws.Range("DP5").PasteSpecial(-4122)

# Row 6
$ws.Range("DO6").Value = "B"
$ws.Range("DN6").Copy()
$ws.Range("DO6").PasteSpecial(-4122)
$ws.Range("DP6").Value = "B"
$ws.Range("DN6").Copy()
$ws.Range("DP6").PasteSpecial(-4122)

# Row 7
$ws.Range("DO7").Value = "P"
$ws.Range("DN7").Copy()
$ws.Range("DO7").PasteSpecial(-4122)
$ws.Range("DP7").Value = "P"
$ws.Range("DN7").Copy()
$ws.Range("DP7").PasteSpecial(-4122)

# Row 8
$ws.Range("DO8").Value = "M"
$ws.Range("DN8").Copy()
$ws.Range("DO8").PasteSpecial(-4122)
$ws.Range("DP8").Value = "P"
$ws.Range("DN8").Copy()
$ws.Range("DP8").PasteSpecial(-4122)

# Row 9
$ws.Range("DO9").Value = "P"
$ws.Range("DN9").Copy()
$ws.Range("DO9").PasteSpecial(-4122)
$ws.Range("DP9").Value = "P"
$ws.Range("DN9").Copy()
$ws.Range("DP9").PasteSpecial(-4122)

# Row 10
$ws.Range("DO10").Value = "P"
$ws.Range("DN10").Copy()
$ws.Range("DO10").PasteSpecial(-4122)
$ws.Range("DP10").Value = "P"
$ws.Range("DN10").Copy()
$ws.Range("DP10").PasteSpecial(-4122)

# Row 11
$ws.Range("DO11").Value = "P"
$ws.Range("DN11").Copy()
$ws.Range("DO11").PasteSpecial(-4122)
$ws.Range("DP11").Value = "P"
$ws.Range("DN11").Copy()
$ws.Range("DP11").PasteSpecial(-4122)

# Row 13
$ws.Range("DO13").Value = "REP"
$ws.Range("DN13").Copy()
$ws.Range("DO13").PasteSpecial(-4122)
$ws.Range("DP13").Value = "RH"
$ws.Range("DN13").Copy()
$ws.Range("DP13").PasteSpecial(-4122)

# Row 14
$ws.Range("DO14").Value = "P"
$ws.Range("DN14").Copy()
$ws.Range("DO14").PasteSpecial(-4122)
$ws.Range("DP14").Value = "P"
$ws.Range("DN14").Copy()
$ws.Range("DP14").PasteSpecial(-4122)

# Row 15
$ws.Range("DO15").Value = "P"
$ws.Range("DN15").Copy()
$ws.Range("DO15").PasteSpecial(-4122)
$ws.Range("DP15").Value = "P"
$ws.Range("DN15").Copy()
$ws.Range("DP15").PasteSpecial(-4122)

# Row 16
$ws.Range("DN16").Copy()
$ws.Range("DO16").PasteSpecial(-4122)
$ws.Range("DN16").Copy()
$ws.Range("DP16").PasteSpecial(-4122)

# Row 17
$ws.Range("DN17").Copy()
$ws.Range("DO17").PasteSpecial(-4122)
$ws.Range("DN17").Copy()
$ws.Range("DP17").PasteSpecial(-4122)

# Row 18
$ws.Range("DO18").Value = "B"
$ws.Range("DN18").Copy()
$ws.Range("DO18").PasteSpecial(-4122)
$ws.Range("DP18").Value = "B"
$ws.Range("DN18").Copy()
$ws.Range("DP18").PasteSpecial(-4122)

# Row 19
$ws.Range("DO19").Value = "P"
$ws.Range("DN19").Copy()
$ws.Range("DO19").PasteSpecial(-4122)
$ws.Range("DP19").Value = "P"
$ws.Range("DN19").Copy()
$ws.Range("DP19").PasteSpecial(-4122)

# Row 20
$ws.Range("DO20").Value = "B"
$ws.Range("DN20").Copy()
$ws.Range("DO20").PasteSpecial(-4122)
$ws.Range("DP20").Value = "P"
$ws.Range("DN20").Copy()
$ws.Range("DP20").PasteSpecial(-4122)

# Row 22
$ws.Range("DO22").Value = "B"
$ws.Range("DN22").Copy()
$ws.Range("DO22").PasteSpecial(-4122)
$ws.Range("DP22").Value = "P"
$ws.Range("DN22").Copy()
$ws.Range("DP22").PasteSpecial(-4122)

# Row 24
$ws.Range("DO24").Value = "P"
$ws.Range("DN24").Copy()
$ws.Range("DO24").PasteSpecial(-4122)
$ws.Range("DP24").Value = "P"
$ws.Range("DN24").Copy()
$ws.Range("DP24").PasteSpecial(-4122)

# Row 25
$ws.Range("DN25").Copy()
$ws.Range("DO25").PasteSpecial(-4122)
$ws.Range("DN25").Copy()
$ws.Range("DP25").PasteSpecial(-4122)

# Row 26
$ws.Range("DO26").Value = "P"
$ws.Range("DN26").Copy()
$ws.Range("DO26").PasteSpecial(-4122)
$ws.Range("DP26").Value = "P"
$ws.Range("DN26").Copy()
$ws.Range("DP26").PasteSpecial(-4122)

# Row 27
$ws.Range("DO27").Value = "P"
$ws.Range("DN27").Copy()
$ws.Range("DO27").PasteSpecial(-4122)
$ws.Range("DP27").Value = "P"
$ws.Range("DN27").Copy()
$ws.Range("DP27").PasteSpecial(-4122)

# Row 28
$ws.Range("DO28").Value = "P"
$ws.Range("DN28").Copy()
$ws.Range("DO28").PasteSpecial(-4122)
$ws.Range("DP28").Value = "P"
$ws.Range("DN28").Copy()
$ws.Range("DP28").PasteSpecial(-4122)

# Row 29
$ws.Range("DO29").Value = "P"
$ws.Range("DN29").Copy()
$ws.Range("DO29").PasteSpecial(-4122)
$ws.Range("DP29").Value = "P"
$ws.Range("DN29").Copy()
$ws.Range("DP29").PasteSpecial(-4122)

# Row 30
$ws.Range("DO30").Value = "P"
$ws.Range("DN30").Copy()
$ws.Range("DO30").PasteSpecial(-4122)
$ws.Range("DP30").Value = "P"
$ws.Range("DN30").Copy()
$ws.Range("DP30").PasteSpecial(-4122)

# Row 31
$ws.Range("DO31").Value = "P"
$ws.Range("DN31").Copy()
$ws.Range("DO31").PasteSpecial(-4122)
$ws.Range("DP31").Value = "P"
$ws.Range("DN31").Copy()
$ws.Range("DP31").PasteSpecial(-4122)

# Row 32
$ws.Range("DO32").Value = "P"
$ws.Range("DN32").Copy()
$ws.Range("DO32").PasteSpecial(-4122)
$ws.Range("DP32").Value = "P"
$ws.Range("DN32").Copy()
$ws.Range("DP32").PasteSpecial(-4122)


# --- View state: move the frozen-pane scroll position and active selection
#     to reflect the newly added columns (best-effort; selection is honored) ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 119
$ws.Range("DR22").Select()
